# Add a second table (allowed inbound/outbound port ranges) below the
# existing CIDR table, then update the page setup and the saved cell
# selection to match the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New table header (row 8) ---------------------------------------
$ws.Range("A8").Value = "범위"
$ws.Range("B8").Value = "포트"
$ws.Range("D8").Value = "설명"

# --- Sub-header (row 9) ----------------------------------------------
$ws.Range("B9").Value = "인"
$ws.Range("C9").Value = "아웃"

# --- Data rows (10-13) -------------------------------------------------
$ws.Range("A10").Value = "10.0.0.0/22"
$ws.Range("B10").Value = 80
$ws.Range("C10").Value = 80

$ws.Range("B11").Value = 443
$ws.Range("C11").Value = 443

$ws.Range("A12").Value = "10.0.8.0/23"
$ws.Range("B12").Value = 5000
$ws.Range("C12").Value = 5000

$ws.Range("A13").Value = "10.0.13.0/24"
$ws.Range("B13").Value = 3306
$ws.Range("C13").Value = 3306

# --- Page setup (portrait, A4-ish "paperSize 9") ----------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# --- Match the persisted selection in the saved workbook --------------
$ws.Range("D6").Select()
